# Vehicle Module completed and assign trailer code commented as UI changes
#
# The "loadingtime" sheet previously tracked four snapshot columns
# ( Sep 16 /  Sep 25 /  Sep 25 /  Sep 26 ). The two trailer-code columns
# are no longer needed now that the Vehicle Module work is done, so they
# are removed and the remaining snapshot column is rolled forward to
# " Sep 27" with its updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailer-code columns (D and E) entirely; remaining columns
# shift left, shrinking the sheet's used range from A1:E4 down to A1:C4.
$ws.Range("D1:E1").EntireColumn.Delete()

# Roll the remaining snapshot column header forward to " Sep 27" and update
# its figures for the new snapshot date.
$ws.Range("C1").Value = " Sep 27"
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 2
